$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$csv = @'
45785,622
45785.01041666666,619
45785.02083333334,620
45785.03125,620
45785.04166666666,621
45785.05208333334,620
45785.0625,621
45785.07291666666,620
45785.08333333334,621
45785.09375,620
45785.10416666666,624
45785.11458333334,622
45785.125,610
45785.13541666666,605
45785.14583333334,605
45785.15625,611
45785.16666666666,625
45785.17708333334,624
45785.1875,624
45785.19791666666,621
45785.20833333334,538
45785.21875,544
45785.22916666666,546
45785.23958333334,573
45785.25,731
45785.26041666666,737
45785.27083333334,734
45785.28125,740
45785.29166666666,775
45785.30208333334,777
45785.3125,785
45785.32291666666,798
45785.33333333334,831
45785.34375,827
45785.35416666666,821
45785.36458333334,799
45785.375,650
45785.38541666666,649
45785.39583333334,648
45785.40625,635
45785.41666666666,459
45785.42708333334,466
45785.4375,442
45785.44791666666,432
45785.45833333334,277
45785.46875,270
45785.47916666666,272
45785.48958333334,266
45785.5,312
45785.51041666666,318
45785.52083333334,328
45785.53125,350
45785.54166666666,340
45785.55208333334,325
45785.5625,302
45785.57291666666,301
45785.58333333334,291
45785.59375,288
45785.60416666666,289
45785.61458333334,289
45785.625,195
45785.63541666666,192
45785.64583333334,227
45785.65625,248
45785.66666666666,434
45785.67708333334,442
45785.6875,462
45785.69791666666,512
45785.70833333334,834
45785.71875,850
45785.72916666666,847
45785.73958333334,839
45785.75,903
45785.76041666666,918
45785.77083333334,914
45785.78125,949
45785.79166666666,908
45785.80208333334,924
45785.8125,933
45785.82291666666,939
45785.83333333334,875
45785.84375,870
45785.85416666666,850
45785.86458333334,867
45785.875,826
45785.88541666666,816
45785.89583333334,775
45785.90625,770
45785.91666666666,516
45785.92708333334,492
45785.9375,470
45785.94791666666,465
45785.95833333334,488
45785.96875,494
45785.97916666666,500
45785.98958333334,585
45786,571
45786.01041666666,581
45786.02083333334,574
45786.03125,556
45786.04166666666,551
45786.05208333334,546
45786.0625,537
45786.07291666666,534
45786.08333333334,548
45786.09375,548
45786.10416666666,589
45786.11458333334,574
45786.125,555
45786.13541666666,556
45786.14583333334,548
45786.15625,553
45786.16666666666,550
45786.17708333334,501
45786.1875,540
45786.19791666666,561
45786.20833333334,341
45786.21875,341
45786.22916666666,345
45786.23958333334,383
45786.25,536
45786.26041666666,0
45786.27083333334,0
45786.28125,0
45786.29166666666,0
45786.30208333334,0
45786.3125,0
45786.32291666666,0
45786.33333333334,0
45786.34375,0
45786.35416666666,0
45786.36458333334,0
45786.375,0
45786.38541666666,0
45786.39583333334,0
45786.40625,0
45786.41666666666,0
45786.42708333334,0
45786.4375,0
45786.44791666666,0
45786.45833333334,0
45786.46875,0
45786.47916666666,0
45786.48958333334,0
45786.5,0
45786.51041666666,0
45786.52083333334,0
45786.53125,0
45786.54166666666,0
45786.55208333334,0
45786.5625,0
45786.57291666666,0
45786.58333333334,0
45786.59375,0
45786.60416666666,0
45786.61458333334,0
45786.625,0
45786.63541666666,0
45786.64583333334,0
45786.65625,0
45786.66666666666,0
45786.67708333334,0
45786.6875,0
45786.69791666666,0
45786.70833333334,0
45786.71875,0
45786.72916666666,0
45786.73958333334,0
45786.75,0
45786.76041666666,0
45786.77083333334,0
45786.78125,0
45786.79166666666,0
45786.80208333334,0
45786.8125,0
45786.82291666666,0
45786.83333333334,0
45786.84375,0
45786.85416666666,0
45786.86458333334,0
45786.875,0
45786.88541666666,0
45786.89583333334,0
45786.90625,0
45786.91666666666,0
45786.92708333334,0
45786.9375,0
45786.94791666666,0
45786.95833333334,0
45786.96875,0
45786.97916666666,0
45786.98958333334,0
'@

$lines = $csv -split "`n"

$row = 2
foreach ($line in $lines) {
    $parts = $line.Trim().Split(',')
    $ws.Cells.Item($row, 1).Value = [double]$parts[0]
    $ws.Cells.Item($row, 2).Value = [double]$parts[1]
    $row++
}
